$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds dates stored as serial numbers.
# All rows from 2 to 181 had their value bumped by one day (45205 -> 45206).
$ws.Range("C2:C181").Value = 45206
